# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.506.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.109.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.90%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "

# Row 6
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5274"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08916"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.112.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.763"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.807"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.65%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001133"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06647"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.325"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.553.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.346"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.355.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.664"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.234"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.704"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +25.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.274"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.903"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02598"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.44%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.526"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2290"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.25%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.25%  "

# Row 43
$ws.Range("E43").Value = "  +3.42%  "

# Row 44
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6443"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.11%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.66%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.229"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.635"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("E49").Value = "  -0.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.219"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.94%  "

